# Auto-generated edit script: updates cached market-price / profit values
# across several sheets (ALC, ARM, BSM, CRP, CUL, LTW, WVR) to reflect a
# scheduled data-refresh run, per the commit 'chore: update Sheets via scheduled runner'.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1750
$ws.Range("I40").Value = 1666.6666
$ws.Range("J40").Value = 1875
$ws.Range("K40").Value = 1666.6666
$ws.Range("L40").Value = 1875
$ws.Range("M40").Value = -1491.6666
$ws.Range("N40").Value = -2225

$ws.Range("H112").Value = 2430.4211
$ws.Range("J112").Value = 2932
$ws.Range("L112").Value = 8796
$ws.Range("N112").Value = -11012

$ws.Range("H132").Value = 2113
$ws.Range("J132").Value = 3500
$ws.Range("L132").Value = 10500
$ws.Range("N132").Value = -15560

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1697.28
$ws.Range("I2").Value = 1594.4706
$ws.Range("K2").Value = 1594.4706
$ws.Range("M2").Value = -1481.4706

$ws.Range("H5").Value = 136.42857
$ws.Range("I5").Value = 196.66667
$ws.Range("J5").Value = 91.25
$ws.Range("K5").Value = 196.66667
$ws.Range("L5").Value = 91.25
$ws.Range("M5").Value = -84.66667000000001
$ws.Range("N5").Value = -315.25

$ws.Range("H116").Value = 1697.28
$ws.Range("I116").Value = 1594.4706
$ws.Range("K116").Value = 1594.4706
$ws.Range("M116").Value = 699.5293999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1697.28
$ws.Range("I3").Value = 1594.4706
$ws.Range("K3").Value = 1594.4706
$ws.Range("M3").Value = -1480.4706

$ws.Range("H4").Value = 136.42857
$ws.Range("I4").Value = 196.66667
$ws.Range("J4").Value = 91.25
$ws.Range("K4").Value = 196.66667
$ws.Range("L4").Value = 91.25
$ws.Range("M4").Value = -81.66667000000001
$ws.Range("N4").Value = -321.25

$ws.Range("H22").Value = 12795.25
$ws.Range("I22").Value = 12795.25
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 12795.25
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -12622.25
$ws.Range("N22").ClearContents()

$ws.Range("H134").Value = 3101.353
$ws.Range("I134").Value = 2824.8462
$ws.Range("K134").Value = 8474.5386
$ws.Range("M134").Value = -5939.5386

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1179.1111
$ws.Range("I16").Value = 1118.5
$ws.Range("J16").Value = 1300.3334
$ws.Range("K16").Value = 1118.5
$ws.Range("L16").Value = 1300.3334
$ws.Range("M16").Value = -831.5
$ws.Range("N16").Value = -1874.3334

$ws.Range("H17").Value = 32333.334
$ws.Range("I17").Value = 15000
$ws.Range("J17").Value = 41000
$ws.Range("K17").Value = 15000
$ws.Range("L17").Value = 41000
$ws.Range("M17").Value = -14826
$ws.Range("N17").Value = -41348

$ws.Range("H22").Value = 394.77777
$ws.Range("I22").Value = 225.875
$ws.Range("J22").Value = 1746
$ws.Range("K22").Value = 225.875
$ws.Range("L22").Value = 1746
$ws.Range("M22").Value = 124.125
$ws.Range("N22").Value = -2446

$ws.Range("H33").Value = 9765.5
$ws.Range("I33").Value = 9765.5
$ws.Range("J33").Value = 0
$ws.Range("K33").Value = 9765.5
$ws.Range("L33").Value = 0
$ws.Range("M33").Value = -9386.5
$ws.Range("N33").ClearContents()

$ws.Range("H107").Value = 426.9375
$ws.Range("I107").Value = 366.5
$ws.Range("J107").Value = 487.375
$ws.Range("K107").Value = 366.5
$ws.Range("L107").Value = 487.375
$ws.Range("M107").Value = 1553.5
$ws.Range("N107").Value = -4327.375

$ws.Range("H113").Value = 1179.1111
$ws.Range("I113").Value = 1118.5
$ws.Range("J113").Value = 1300.3334
$ws.Range("K113").Value = 1118.5
$ws.Range("L113").Value = 1300.3334
$ws.Range("M113").Value = 1051.5
$ws.Range("N113").Value = -5640.3334

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 638.7843
$ws.Range("I113").Value = 596.8484999999999
$ws.Range("J113").Value = 715.6667
$ws.Range("K113").Value = 1790.5455
$ws.Range("L113").Value = 2147.0001
$ws.Range("M113").Value = 379.4545000000003
$ws.Range("N113").Value = -6487.0001

$ws.Range("H122").Value = 286
$ws.Range("I122").Value = 290.4
$ws.Range("J122").Value = 275
$ws.Range("K122").Value = 2613.6
$ws.Range("L122").Value = 2475
$ws.Range("M122").Value = -163.5999999999999
$ws.Range("N122").Value = -7375

$ws.Range("H136").Value = 6077.8667
$ws.Range("I136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("M136").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 9432
$ws.Range("I2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("M2").ClearContents()

$ws.Range("H22").Value = 1200.3334
$ws.Range("I22").Value = 2001
$ws.Range("J22").Value = 800
$ws.Range("K22").Value = 2001
$ws.Range("L22").Value = 800
$ws.Range("M22").Value = -1706
$ws.Range("N22").Value = -1390

$ws.Range("H27").Value = 1200.3334
$ws.Range("I27").Value = 2001
$ws.Range("J27").Value = 800
$ws.Range("K27").Value = 2001
$ws.Range("L27").Value = 800
$ws.Range("M27").Value = -1894
$ws.Range("N27").Value = -1014

$ws.Range("H61").Value = 28513.5
$ws.Range("I61").Value = 36701.332
$ws.Range("J61").Value = 3950
$ws.Range("K61").Value = 36701.332
$ws.Range("L61").Value = 3950
$ws.Range("M61").Value = -36499.332
$ws.Range("N61").Value = -4354

$ws.Range("H93").Value = 2694.7
$ws.Range("I93").Value = 2609
$ws.Range("K93").Value = 2609
$ws.Range("M93").Value = -1361

$ws.Range("H113").Value = 28513.5
$ws.Range("I113").Value = 36701.332
$ws.Range("J113").Value = 3950
$ws.Range("K113").Value = 36701.332
$ws.Range("L113").Value = 3950
$ws.Range("M113").Value = -34531.332
$ws.Range("N113").Value = -8290

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 220760
$ws.Range("I81").Value = 220760
$ws.Range("K81").Value = 441520
$ws.Range("M81").Value = -440459

$ws.Range("H84").Value = 220760
$ws.Range("I84").Value = 220760
$ws.Range("K84").Value = 2207600
$ws.Range("M84").Value = -2202296

$ws.Range("H113").Value = 679.1875
$ws.Range("I113").Value = 629.1111
$ws.Range("J113").Value = 743.5714
$ws.Range("K113").Value = 1887.3333
$ws.Range("L113").Value = 2230.7142
$ws.Range("M113").Value = 282.6667000000002
$ws.Range("N113").Value = -6570.7142

$ws.Range("H132").Value = 2445.818
$ws.Range("I132").Value = 1767.7222
$ws.Range("J132").Value = 3259.5334
$ws.Range("K132").Value = 5303.1666
$ws.Range("L132").Value = 9778.600199999999
$ws.Range("M132").Value = -2773.1666
$ws.Range("N132").Value = -14838.6002
